# Refactored shared assertions in SitewideSearch_Test
#
# The "argle-bargle" expected value in the SitewideSearch sheet is replaced
# with the broader "argle-bargle or foofaraw" assertion text.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SitewideSearch")
$ws.Range("E2").Value = "argle-bargle or foofaraw"
